# Code_Verification_Status.xlsx update:
# - "XX" sheet: replace the 4b/4c/4d_run_mSigHdp.R rows with the new
#   4c_run_signeR.R / 4b_run_SignatureAnalyzer.R / 4d_run_SP.py rows,
#   add a "branch" comment next to 4a_run_mSigHdp.R, mark the newly
#   touched status cells with the "fixed" (orange) color, and drop the
#   purple highlight from the footnote-marker cells.
# - Update the active selection on the "XX" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("XX")

# Orange "Fixed errors in path and comments" fill color (matches style used
# elsewhere in the sheet for the same status, e.g. A10 on the README sheet).
$orange = 49407

# Row 14: 4a_run_mSigHdp.R keeps its name, gains a comment about branches.
$ws.Range("B14").Value = 'Need to change branch to use as "master"'

# Row 15 used to be 4b_run_mSigHdp.R; it becomes the "*" footnote marker row,
# with its status cells colored to show "fixed errors" (not yet verified).
$ws.Range("A15").Value = "*"
$ws.Range("B15:E15").Interior.Color = $orange

# Row 16 used to be 4c_run_mSigHdp.R; it becomes 4c_run_signeR.R.
$ws.Range("A16").Value = "4c_run_signeR.R"
$ws.Range("B16:E16").Interior.Color = $orange

# Row 17 used to be 4d_run_mSigHdp.R; it becomes 4d_run_SP.py.
$ws.Range("A17").Value = "4d_run_SP.py"
$ws.Range("B17:E17").Interior.Color = $orange

# Row 18 (4e_SBS_NR_hdp_gamma_beta_20.R): B18:C18 newly colored orange,
# D18:E18 keep their existing "N/A" content/format.
$ws.Range("B18:C18").Interior.Color = $orange

# Row 19 (4e_SBS_NR_hdp_gamma_beta_50.R): D19:E19 newly colored orange,
# B19:C19 keep their existing "N/A" content/format.
$ws.Range("D19:E19").Interior.Color = $orange

# Row 20 (4f_SBS_NR_hdp_gamma_beta_1.R): B20:E20 all become empty/orange
# (previously B20:C20 held "N/A").
$ws.Range("B20:C20").ClearContents()
$ws.Range("B20:E20").Interior.Color = $orange

# Row 21 (5_rename_SA_SP_files.R): B21:E21 all become empty/orange.
$ws.Range("B21:E21").Interior.Color = $orange

# Row 11: the "*" footnote-marker cells (B11/D11) drop their purple fill,
# keeping the red font.
$ws.Range("B11").Interior.Pattern = -4142
$ws.Range("D11").Interior.Pattern = -4142

# Update the selection / scroll position on the "XX" sheet.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("E21").Select()
